# "committing this so I can switch to a new branch"
# Adds a 4th data column (J) to the Sprint-3 retro table (rows 79:91),
# appends the Sprint 4 planning table (rows 93:127) on Sheet2, extends
# Chart 4's series to include the new column, and widens column F to
# fit the longer sprint-4 task labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$xlLeft = -4131

# ---------------------------------------------------------------
# 1. Sprint 3 retro table: add the "week 3" column (J) to rows 80-91
# ---------------------------------------------------------------
$ws.Range("J80").Value = 1
$ws.Range("J81").Value = 1
$ws.Range("J82").Value = 8
$ws.Range("J83").Value = 0
$ws.Range("J84").Value = 5
$ws.Range("J85").Value = 1
$ws.Range("J86").Value = 1
$ws.Range("J87").Value = 15
$ws.Range("J88").Value = 4
$ws.Range("J89").Value = 1
$ws.Range("J90").Value = 2
$ws.Range("J91").Formula = "=SUM(J80:J90)"

# ---------------------------------------------------------------
# 2. New Sprint 4 planning table, rows 93-127
# ---------------------------------------------------------------
$ws.Range("F93").Value = "Sprint 4"

$ws.Range("F94").Value = "Week"
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 1
$ws.Range("I94").Value = 2
$ws.Range("J94").Value = 3

$ws.Range("F95").Value = "Karly"

$ws.Range("F96").Value = "Change resolution, switch to android"
$ws.Range("G96").Value = 2

$ws.Range("F97").Value = "Checksum game termination"
$ws.Range("G97").Value = 1

$ws.Range("F98").Value = "Checksum level integration"
$ws.Range("G98").Value = 1

$ws.Range("F99").Value = "Dijkstra testing"
$ws.Range("G99").Value = 4

$ws.Range("F100").Value = "Segmentation Backend"
$ws.Range("G100").Value = 4

$ws.Range("F101").Value = "dialog writing"

$ws.Range("F102").Value = "checksum explanation"
$ws.Range("F102").HorizontalAlignment = $xlLeft
$ws.Range("F102").IndentLevel = 1
$ws.Range("G102").Value = 1

$ws.Range("F103").Value = "binary addition instruction"
$ws.Range("F103").HorizontalAlignment = $xlLeft
$ws.Range("F103").IndentLevel = 1
$ws.Range("G103").Value = 1

$ws.Range("F104").Value = "pop-up encouragement"
$ws.Range("F104").HorizontalAlignment = $xlLeft
$ws.Range("F104").IndentLevel = 1
$ws.Range("G104").Value = 1

$ws.Range("F105").Value = "Dijkstra explanation"
$ws.Range("F105").HorizontalAlignment = $xlLeft
$ws.Range("F105").IndentLevel = 1
$ws.Range("G105").Value = 1

$ws.Range("F106").Value = "pathfind instruction"
$ws.Range("F106").HorizontalAlignment = $xlLeft
$ws.Range("F106").IndentLevel = 1
$ws.Range("G106").Value = 1

$ws.Range("F107").Value = "Segmentation explanation"
$ws.Range("F107").HorizontalAlignment = $xlLeft
$ws.Range("F107").IndentLevel = 1
$ws.Range("G107").Value = 1

$ws.Range("F108").Value = "package instructions"
$ws.Range("F108").HorizontalAlignment = $xlLeft
$ws.Range("F108").IndentLevel = 1
$ws.Range("G108").Value = 1

$ws.Range("F109").Value = "Final background design"
$ws.Range("F109").HorizontalAlignment = $xlLeft
$ws.Range("G109").Value = 2

$ws.Range("F110").Value = "Dijkstra scoring"
$ws.Range("F110").HorizontalAlignment = $xlLeft
$ws.Range("G110").Value = 3

$ws.Range("F111").Value = "Segmentation scoring"
$ws.Range("F111").HorizontalAlignment = $xlLeft
$ws.Range("G111").Value = 3

$ws.Range("F113").Value = "Nef"

$ws.Range("F114").Value = "pathfind front end"

$ws.Range("F115").Value = "add edge bubbles"
$ws.Range("F115").HorizontalAlignment = $xlLeft
$ws.Range("F115").IndentLevel = 1
$ws.Range("G115").Value = 0.5

$ws.Range("F116").Value = "renumber edges"
$ws.Range("F116").HorizontalAlignment = $xlLeft
$ws.Range("F116").IndentLevel = 1
$ws.Range("G116").Value = 0.5

$ws.Range("F117").Value = "tie edges to backend"
$ws.Range("F117").HorizontalAlignment = $xlLeft
$ws.Range("F117").IndentLevel = 1
$ws.Range("G117").Value = 0.5

$ws.Range("F118").Value = "Segmentation front end"
$ws.Range("F118").HorizontalAlignment = $xlLeft

$ws.Range("F119").Value = "scrolling boxes"
$ws.Range("F119").HorizontalAlignment = $xlLeft
$ws.Range("F119").IndentLevel = 1
$ws.Range("G119").Value = 2

$ws.Range("F120").Value = "click and drag boxes"
$ws.Range("F120").HorizontalAlignment = $xlLeft
$ws.Range("F120").IndentLevel = 1
$ws.Range("G120").Value = 4

$ws.Range("F121").Value = "trash can image changing"
$ws.Range("F121").HorizontalAlignment = $xlLeft
$ws.Range("F121").IndentLevel = 1
$ws.Range("G121").Value = 2

$ws.Range("F122").Value = "score box updating"
$ws.Range("F122").HorizontalAlignment = $xlLeft
$ws.Range("F122").IndentLevel = 1
$ws.Range("G122").Value = 1

$ws.Range("F123").Value = "checksum front end move numbers"
$ws.Range("F123").HorizontalAlignment = $xlLeft
$ws.Range("G123").Value = 1

$ws.Range("F124").Value = "Checksum testing"
$ws.Range("F124").HorizontalAlignment = $xlLeft
$ws.Range("G124").Value = 4

$ws.Range("G127").Formula = "=SUM(G96:G124)"

# ---------------------------------------------------------------
# 3. Widen column F so the longer sprint-4 labels fit
# ---------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 28.6666666666667

# ---------------------------------------------------------------
# 4. Extend Chart 4's series so it plots the new "week 3" column
# ---------------------------------------------------------------
$chart4 = $ws.ChartObjects().Item(4).Chart
$series = $chart4.SeriesCollection().Item(1)
$series.Formula = "=SERIES(,Sheet2!`$G`$79:`$J`$79,Sheet2!`$G`$91:`$J`$91,1)"

# ---------------------------------------------------------------
# 5. Leave the selection where the author left it
# ---------------------------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("G112").Select()
